$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "41.947.46"
$ws.Range("E2").Value = "  -2.17%  "

# Row 3
$ws.Range("D3").Value = "2.296.82"
$ws.Range("E3").Value = "  -2.47%  "

# Row 4
$ws.Range("E4").Value = "  -0.04%  "

# Row 5
$ws.Range("D5").Value = "'315.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.07%  "

# Row 6
$ws.Range("D6").Value = "'103.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.71%  "

# Row 7
$ws.Range("D7").Value = "'0.626"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.95%  "

# Row 8
$ws.Range("E8").Value = "  +0.04%  "

# Row 9
$ws.Range("D9").Value = "'0.608"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.69%  "

# Row 10
$ws.Range("D10").Value = "'39.90"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.15%  "

# Row 11
$ws.Range("D11").Value = "'0.0909"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.80%  "

# Row 12
$ws.Range("D12").Value = "'8.41"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.98%  "

# Row 13
$ws.Range("E13").Value = "  +0.90%  "

# Row 14
$ws.Range("E14").Value = "  -2.49%  "

# Row 15
$ws.Range("D15").Value = "'15.33"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.28%  "

# Row 16
$ws.Range("D16").Value = "2.642.99"
$ws.Range("E16").Value = "  -2.56%  "

# Row 17
$ws.Range("D17").Value = "2.300.14"
$ws.Range("E17").Value = "  -1.37%  "

# Row 18
$ws.Range("D18").Value = "41.920.69"
$ws.Range("E18").Value = "  -2.19%  "

# Row 19
$ws.Range("D19").Value = "'7.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.65%  "

# Row 20
$ws.Range("E20").Value = "  -1.05%  "

# Row 21
$ws.Range("D21").Value = "'72.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.92%  "

# Row 22
$ws.Range("E22").Value = "  -2.41%  "

# Row 23
$ws.Range("D23").Value = "'257.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.12%  "

# Row 24
$ws.Range("E24").Value = "  -0.64%  "

# Row 25
$ws.Range("D25").Value = "'9.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.26%  "

# Row 26
$ws.Range("E26").Value = "  +0.68%  "

# Row 27
$ws.Range("D27").Value = "'10.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.18%  "

# Row 28
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").Value = "'2.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.67%  "

# Row 29
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'22.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.26%  "

# Row 30
$ws.Range("D30").Value = "'35.54"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.38%  "

# Row 31
$ws.Range("D31").Value = "'164.51"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.76%  "

# Row 32
$ws.Range("D32").Value = "'0.0886"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.62%  "

# Row 33
$ws.Range("E33").Value = "  -3.04%  "

# Row 34
$ws.Range("D34").Value = "'5.86"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.26%  "

# Row 35
$ws.Range("E35").Value = "  -0.54%  "

# Row 36
$ws.Range("D36").Value = "'0.118"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.09%  "

# Row 37
$ws.Range("E37").Value = "  -0.26%  "

# Row 38
$ws.Range("E38").Value = "  +8.04%  "

# Row 39
$ws.Range("E39").Value = "  -2.96%  "

# Row 40
$ws.Range("D40").Value = "'3.66"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.92%  "

# Row 41
$ws.Range("D41").Value = "'100.56"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +17.08%  "

# Row 42
$ws.Range("E42").Value = "  +1.00%  "

# Row 43
$ws.Range("D43").Value = "'70.20"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.16%  "

# Row 44
$ws.Range("D44").Value = "'0.226"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.59%  "

# Row 45
$ws.Range("E45").Value = "  -0.06%  "

# Row 46
$ws.Range("D46").Value = "'12.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.87%  "

# Row 47
$ws.Range("D47").Value = "'114.22"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.80%  "

# Row 48
$ws.Range("D48").Value = "'77.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +5.18%  "

# Row 49
$ws.Range("D49").Value = "'9.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.81%  "

# Row 50
$ws.Range("E50").Value = "  -4.06%  "

# Row 51
$ws.Range("E51").Value = "  +2.15%  "
